$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price column (D) values, preserving text formatting ---
$style_D2 = $ws.Range("D2").Style
$style_D3 = $ws.Range("D3").Style
$style_D5 = $ws.Range("D5").Style
$style_D7 = $ws.Range("D7").Style
$style_D8 = $ws.Range("D8").Style
$style_D9 = $ws.Range("D9").Style
$style_D10 = $ws.Range("D10").Style
$style_D11 = $ws.Range("D11").Style
$style_D12 = $ws.Range("D12").Style
$style_D13 = $ws.Range("D13").Style
$style_D14 = $ws.Range("D14").Style
$style_D16 = $ws.Range("D16").Style
$style_D17 = $ws.Range("D17").Style
$style_D18 = $ws.Range("D18").Style
$style_D19 = $ws.Range("D19").Style
$style_D20 = $ws.Range("D20").Style
$style_D22 = $ws.Range("D22").Style
$style_D23 = $ws.Range("D23").Style
$style_D24 = $ws.Range("D24").Style
$style_D25 = $ws.Range("D25").Style
$style_D26 = $ws.Range("D26").Style
$style_D27 = $ws.Range("D27").Style
$style_D29 = $ws.Range("D29").Style
$style_D30 = $ws.Range("D30").Style
$style_D32 = $ws.Range("D32").Style
$style_D33 = $ws.Range("D33").Style
$style_D34 = $ws.Range("D34").Style
$style_D36 = $ws.Range("D36").Style
$style_D38 = $ws.Range("D38").Style
$style_D41 = $ws.Range("D41").Style
$style_D42 = $ws.Range("D42").Style
$style_D43 = $ws.Range("D43").Style
$style_D44 = $ws.Range("D44").Style
$style_D45 = $ws.Range("D45").Style
$style_D47 = $ws.Range("D47").Style
$style_D48 = $ws.Range("D48").Style
$style_D49 = $ws.Range("D49").Style
$style_D50 = $ws.Range("D50").Style
$style_D51 = $ws.Range("D51").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.512.63"
$ws.Range("D3").Value = "1.826.76"
$ws.Range("D5").Value = "316.62"
$ws.Range("D7").Value = "0.5409"
$ws.Range("D8").Value = "0.4045"
$ws.Range("D9").Value = "0.07649"
$ws.Range("D10").Value = "1.120"
$ws.Range("D11").Value = "41.90"
$ws.Range("D12").Value = "6.328"
$ws.Range("D13").Value = "7.647"
$ws.Range("D14").Value = "20.99"
$ws.Range("D16").Value = "1.823.88"
$ws.Range("D17").Value = "0.00001089"
$ws.Range("D18").Value = "89.95"
$ws.Range("D19").Value = "0.06605"
$ws.Range("D20").Value = "17.71"
$ws.Range("D22").Value = "6.074"
$ws.Range("D23").Value = "28.511.02"
$ws.Range("D24").Value = "11.15"
$ws.Range("D25").Value = "2.271"
$ws.Range("D26").Value = "157.79"
$ws.Range("D27").Value = "2.452"
$ws.Range("D29").Value = "2.034.22"
$ws.Range("D30").Value = "124.02"
$ws.Range("D32").Value = "0.1106"
$ws.Range("D33").Value = "5.677"
$ws.Range("D34").Value = "0.07434"
$ws.Range("D36").Value = "0.2236"
$ws.Range("D38").Value = "5.215"
$ws.Range("D41").Value = "0.6293"
$ws.Range("D42").Value = "1.186"
$ws.Range("D43").Value = "0.9997"
$ws.Range("D44").Value = "1.397"
$ws.Range("D45").Value = "13.47"
$ws.Range("D47").Value = "0.5865"
$ws.Range("D48").Value = "125.40"
$ws.Range("D49").Value = "2.005"
$ws.Range("D50").Value = "1.197"
$ws.Range("D51").Value = "0.06885"

$ws.Range("D2").Style = $style_D2
$ws.Range("D3").Style = $style_D3
$ws.Range("D5").Style = $style_D5
$ws.Range("D7").Style = $style_D7
$ws.Range("D8").Style = $style_D8
$ws.Range("D9").Style = $style_D9
$ws.Range("D10").Style = $style_D10
$ws.Range("D11").Style = $style_D11
$ws.Range("D12").Style = $style_D12
$ws.Range("D13").Style = $style_D13
$ws.Range("D14").Style = $style_D14
$ws.Range("D16").Style = $style_D16
$ws.Range("D17").Style = $style_D17
$ws.Range("D18").Style = $style_D18
$ws.Range("D19").Style = $style_D19
$ws.Range("D20").Style = $style_D20
$ws.Range("D22").Style = $style_D22
$ws.Range("D23").Style = $style_D23
$ws.Range("D24").Style = $style_D24
$ws.Range("D25").Style = $style_D25
$ws.Range("D26").Style = $style_D26
$ws.Range("D27").Style = $style_D27
$ws.Range("D29").Style = $style_D29
$ws.Range("D30").Style = $style_D30
$ws.Range("D32").Style = $style_D32
$ws.Range("D33").Style = $style_D33
$ws.Range("D34").Style = $style_D34
$ws.Range("D36").Style = $style_D36
$ws.Range("D38").Style = $style_D38
$ws.Range("D41").Style = $style_D41
$ws.Range("D42").Style = $style_D42
$ws.Range("D43").Style = $style_D43
$ws.Range("D44").Style = $style_D44
$ws.Range("D45").Style = $style_D45
$ws.Range("D47").Style = $style_D47
$ws.Range("D48").Style = $style_D48
$ws.Range("D49").Style = $style_D49
$ws.Range("D50").Style = $style_D50
$ws.Range("D51").Style = $style_D51

# --- Update Volume(1h) column (E) values ---
$ws.Range("E2").Value = "  +1.53%  "
$ws.Range("E3").Value = "  +1.87%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  +7.33%  "
$ws.Range("E9").Value = "  +2.76%  "
$ws.Range("E10").Value = "  +2.54%  "
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("E12").Value = "  +3.77%  "
$ws.Range("E13").Value = "  +5.97%  "
$ws.Range("E14").Value = "  +1.95%  "
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("E16").Value = "  +1.83%  "
$ws.Range("E17").Value = "  +3.09%  "
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("E19").Value = "  +2.05%  "
$ws.Range("E20").Value = "  +2.66%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("E22").Value = "  +3.02%  "
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("E25").Value = "  +8.14%  "
$ws.Range("E26").Value = "  +1.88%  "
$ws.Range("E27").Value = "  +7.43%  "
$ws.Range("E28").Value = "  +2.61%  "
$ws.Range("E29").Value = "  +2.11%  "
$ws.Range("E30").Value = "  +2.87%  "
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("E32").Value = "  +4.86%  "
$ws.Range("E33").Value = "  +2.42%  "
$ws.Range("E34").Value = "  +14.17%  "
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("E37").Value = "  +2.82%  "
$ws.Range("E38").Value = "  +4.08%  "
$ws.Range("E39").Value = "  +4.87%  "
$ws.Range("E40").Value = "  +2.88%  "
$ws.Range("E41").Value = "  +2.13%  "
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("E44").Value = "  -3.44%  "
$ws.Range("E45").Value = "  +1.38%  "
$ws.Range("E46").Value = "  +0.73%  "
$ws.Range("E47").Value = "  +1.55%  "
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("E49").Value = "  +4.59%  "
$ws.Range("E50").Value = "  +0.86%  "
$ws.Range("E51").Value = "  +1.06%  "
